# feat: add 2022-Q3 data
#
# The existing "2022-Q2" sheet is duplicated (so the sheetId/relationship-id
# bookkeeping matches what Excel produces when you right-click > Move or
# Copy > Create a copy): the ORIGINAL sheet object keeps its relationship id
# and is turned into the new "2022-Q3" sheet (new quarter data + header
# styling borrowed from the "总计" sheet), while the freshly created COPY
# keeps the name "2022-Q2" and its original data/formatting untouched.
#
# A new row is also inserted into the "总计" (grand-total) sheet holding the
# 2022-Q3 totals, pushing the old 2022-Q2 total row down to row 3.

$wb = $excel.ActiveWorkbook

# xlPasteFormats
$xlPasteFormats = -4122

# Helper: write a value into a cell while forcing it to be stored as TEXT
# (not auto-converted to a number), and without leaving behind the
# cell's previous formatting. Numeric-looking strings such as fund codes
# ("014273") or decimal figures ("3.37") must stay literal text.
function Set-TextValue {
    param($range, $text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# --- duplicate the existing "2022-Q2" sheet, placing the copy right after it ---
$oldQ2 = $wb.Worksheets.Item("2022-Q2")
$oldQ2.Copy($null, $oldQ2)

# the newly created copy is named "2022-Q2 (2)"; rename sheets so that the
# original becomes "2022-Q3" (filled with the new quarter's data) and the
# copy keeps the name "2022-Q2" (its data/format is left exactly as it was).
$newQ2 = $wb.Worksheets.Item("2022-Q2 (2)")
$oldQ2.Name = "2022-Q3"
$newQ2.Name = "2022-Q2"

$q3 = $oldQ2
$total = $wb.Worksheets.Item("总计")

# "2022-Q3" should look like "总计" (page margins + header/row-label style),
# not like the template "2022-Q2" sheet it was duplicated from.
$q3.PageSetup.LeftMargin = $total.PageSetup.LeftMargin
$q3.PageSetup.RightMargin = $total.PageSetup.RightMargin
$q3.PageSetup.TopMargin = $total.PageSetup.TopMargin
$q3.PageSetup.BottomMargin = $total.PageSetup.BottomMargin
$q3.PageSetup.HeaderMargin = $total.PageSetup.HeaderMargin
$q3.PageSetup.FooterMargin = $total.PageSetup.FooterMargin

# --- update the "2022-Q3" sheet with the new quarter's fund data ---

# header row
$total.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial($xlPasteFormats)
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# row 2
$total.Range("A2").Copy()
$q3.Range("A2").PasteSpecial($xlPasteFormats)
$q3.Range("A2").Value = 0
Set-TextValue $q3.Range("B2") "014273"
$q3.Range("C2").Value = "广发北交所精选两年定开混合A"
Set-TextValue $q3.Range("D2") "3.37"
Set-TextValue $q3.Range("E2") "64.25"
Set-TextValue $q3.Range("F2") "3.49"
Set-TextValue $q3.Range("G2") "0.1176"
$q3.Range("H2").Value = 9

# row 3
$total.Range("A2").Copy()
$q3.Range("A3").PasteSpecial($xlPasteFormats)
$q3.Range("A3").Value = 1
Set-TextValue $q3.Range("B3") "014274"
$q3.Range("C3").Value = "广发北交所精选两年定开混合C"
Set-TextValue $q3.Range("D3") "0.85"
Set-TextValue $q3.Range("E3") "64.25"
Set-TextValue $q3.Range("F3") "3.49"
Set-TextValue $q3.Range("G3") "0.0297"
$q3.Range("H3").Value = 9

# --- update the "总计" sheet: row 2 becomes 2022-Q3, row 3 (new) gets the
#     old 2022-Q2 total figures that used to live in row 2 ---

$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial($xlPasteFormats)
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.25

$total.Range("B2").Value = "2022-Q3"
$total.Range("D2").Value = 0.15
